$wb = $excel.ActiveWorkbook

# Sheet ALC, row 32 (item 5484)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1730
$ws.Range("J32").Value = 1625
$ws.Range("L32").Value = 1625
$ws.Range("N32").Value = -2277

# Sheet ALC, row 53 (item 5479)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 554.5
$ws.Range("I53").Value = 243.66667
$ws.Range("J53").Value = 1114
$ws.Range("K53").Value = 243.66667
$ws.Range("L53").Value = 1114
$ws.Range("M53").Value = 393.33333
$ws.Range("N53").Value = -2388

# Sheet ALC, row 138 (item 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 10543.131
$ws.Range("J138").Value = 10879.436
$ws.Range("L138").Value = 32638.308
$ws.Range("N138").Value = -42918.308

# Sheet ARM, row 2 (item 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1126.7142
$ws.Range("I2").Value = 1126.7142
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1126.7142
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1013.7142
$ws.Range("N2").ClearContents()

# Sheet ARM, row 102 (item 19945)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1324.375
$ws.Range("I102").Value = 1324.375
$ws.Range("K102").Value = 1324.375
$ws.Range("M102").Value = 297.625

# Sheet ARM, row 116 (item 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1126.7142
$ws.Range("I116").Value = 1126.7142
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1126.7142
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1167.2858
$ws.Range("N116").ClearContents()

# Sheet ARM, row 132 (item 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3082.7144
$ws.Range("I132").Value = 2346.5
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 7039.5
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -4509.5
$ws.Range("N132").Value = -27560

# Sheet BSM, row 3 (item 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1126.7142
$ws.Range("I3").Value = 1126.7142
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1126.7142
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1012.7142
$ws.Range("N3").ClearContents()

# Sheet BSM, row 64 (item 14184)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 234.66667
$ws.Range("I64").Value = 198
$ws.Range("K64").Value = 198
$ws.Range("M64").Value = 27

# Sheet BSM, row 67 (item 14184)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 234.66667
$ws.Range("I67").Value = 198
$ws.Range("K67").Value = 198
$ws.Range("M67").Value = 582

# Sheet BSM, row 94 (item 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1848.6842
$ws.Range("I94").Value = 1487.7333
$ws.Range("J94").Value = 3202.25
$ws.Range("K94").Value = 1487.7333
$ws.Range("L94").Value = 3202.25
$ws.Range("M94").Value = -1036.7333
$ws.Range("N94").Value = -4104.25

# Sheet CRP, row 16 (item 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2497.1667
$ws.Range("I16").Value = 831.3333
$ws.Range("J16").Value = 4163
$ws.Range("K16").Value = 831.3333
$ws.Range("L16").Value = 4163
$ws.Range("M16").Value = -544.3333
$ws.Range("N16").Value = -4737

# Sheet CRP, row 31 (item 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2902.1177
$ws.Range("I31").Value = 2799.1428
$ws.Range("J31").Value = 3382.6667
$ws.Range("K31").Value = 2799.1428
$ws.Range("L31").Value = 3382.6667
$ws.Range("M31").Value = -2504.1428
$ws.Range("N31").Value = -3972.6667

# Sheet CRP, row 34 (item 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2902.1177
$ws.Range("I34").Value = 2799.1428
$ws.Range("J34").Value = 3382.6667
$ws.Range("K34").Value = 2799.1428
$ws.Range("L34").Value = 3382.6667
$ws.Range("M34").Value = -2597.1428
$ws.Range("N34").Value = -3786.6667

# Sheet CRP, row 113 (item 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2497.1667
$ws.Range("I113").Value = 831.3333
$ws.Range("J113").Value = 4163
$ws.Range("K113").Value = 831.3333
$ws.Range("L113").Value = 4163
$ws.Range("M113").Value = 1338.6667
$ws.Range("N113").Value = -8503

# Sheet CRP, row 132 (item 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7139.64
$ws.Range("J132").Value = 11660.167
$ws.Range("L132").Value = 34980.501
$ws.Range("N132").Value = -40040.501

# Sheet CUL, row 4 (item 4650)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 147.6
$ws.Range("I4").Value = 147.6
$ws.Range("K4").Value = 442.8
$ws.Range("M4").Value = -330.8

# Sheet CUL, row 5 (item 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 610.5
$ws.Range("I5").Value = 444
$ws.Range("K5").Value = 1332
$ws.Range("M5").Value = -1220

# Sheet CUL, row 7 (item 4728)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 472.25
$ws.Range("I7").Value = 94.5
$ws.Range("J7").Value = 850
$ws.Range("K7").Value = 283.5
$ws.Range("L7").Value = 2550
$ws.Range("M7").Value = -171.5
$ws.Range("N7").Value = -2774

# Sheet CUL, row 135 (item 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 610.5
$ws.Range("I135").Value = 444
$ws.Range("K135").Value = 3996
$ws.Range("M135").Value = -1461

# Sheet GSM, row 34 (item 10924)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 37783.5
$ws.Range("J34").Value = 37783.5
$ws.Range("L34").Value = 37783.5
$ws.Range("N34").Value = -38319.5

# Sheet GSM, row 76 (item 10924)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H76").Value = 37783.5
$ws.Range("J76").Value = 37783.5
$ws.Range("L76").Value = 37783.5
$ws.Range("N76").Value = -38413.5

# Sheet GSM, row 79 (item 10924)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H79").Value = 37783.5
$ws.Range("J79").Value = 37783.5
$ws.Range("L79").Value = 37783.5
$ws.Range("N79").Value = -39967.5

# Sheet GSM, row 126 (item 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 224694060
$ws.Range("I126").Value = 224694060
$ws.Range("K126").Value = 674082180
$ws.Range("M126").Value = -674079710

# Sheet GSM, row 132 (item 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4101.25
$ws.Range("I132").Value = 2468.5
$ws.Range("J132").Value = 8999.5
$ws.Range("K132").Value = 7405.5
$ws.Range("L132").Value = 26998.5
$ws.Range("M132").Value = -4875.5
$ws.Range("N132").Value = -32058.5

# Sheet LTW, row 16 (item 5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1123.2
$ws.Range("I16").Value = 1197.091
$ws.Range("J16").Value = 920
$ws.Range("K16").Value = 1197.091
$ws.Range("L16").Value = 920
$ws.Range("M16").Value = -1027.091
$ws.Range("N16").Value = -1260

# Sheet LTW, row 22 (item 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2430.7144
$ws.Range("I22").Value = 746.25
$ws.Range("K22").Value = 746.25
$ws.Range("M22").Value = -451.25

# Sheet LTW, row 27 (item 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2430.7144
$ws.Range("I27").Value = 746.25
$ws.Range("K27").Value = 746.25
$ws.Range("M27").Value = -639.25

# Sheet LTW, row 61 (item 27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# Sheet LTW, row 104 (item 18675)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Sheet LTW, row 113 (item 27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# Sheet LTW, row 136 (item 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4259.9165
$ws.Range("I136").Value = 3902.111
$ws.Range("K136").Value = 11706.333
$ws.Range("M136").Value = -9156.332999999999

# Sheet WVR, row 113 (item 27752)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1181.1111
$ws.Range("I113").Value = 795
$ws.Range("J113").Value = 1490
$ws.Range("K113").Value = 2385
$ws.Range("L113").Value = 4470
$ws.Range("M113").Value = -215
$ws.Range("N113").Value = -8810

# Sheet WVR, row 126 (item 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2167.5
$ws.Range("I126").Value = 1709.5
$ws.Range("J126").Value = 3999.5
$ws.Range("K126").Value = 5128.5
$ws.Range("L126").Value = 11998.5
$ws.Range("M126").Value = -2658.5
$ws.Range("N126").Value = -16938.5

# Sheet WVR, row 132 (item 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3425.6
$ws.Range("J132").Value = 2775
$ws.Range("L132").Value = 8325
$ws.Range("N132").Value = -13385

# Sheet WVR, row 136 (item 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 28892.053
$ws.Range("I136").Value = 29941.611
$ws.Range("K136").Value = 89824.833
$ws.Range("M136").Value = -87274.833
